# Refresh the NATMI ligand-receptor TPM-derived values for the Wnt6-Fzd7
# sheet (script was re-run with new TPM input). Sending/receptor cluster
# labels (columns A-D) are unchanged; only the recalculated expression /
# specificity metrics (columns G-J, M-T) are updated to their new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1693853333333334
$ws.Range("H2").Value = 0.508156
$ws.Range("I2").Value = 0.936523909045002
$ws.Range("J2").Value = 0.936523909045002
$ws.Range("M2").Value = 0.7339303333333334
$ws.Range("N2").Value = 2.201791
$ws.Range("O2").Value = 0.03574007706012852
$ws.Range("P2").Value = 0.03574007706012852
$ws.Range("Q2").Value = 0.1243170341551111
$ws.Range("R2").Value = 1.118853307396
$ws.Range("S2").Value = 0.03347143667792117
$ws.Range("T2").Value = 0.03347143667792116

# Row 3
$ws.Range("G3").Value = 0.1693853333333334
$ws.Range("H3").Value = 0.508156
$ws.Range("I3").Value = 0.936523909045002
$ws.Range("J3").Value = 0.936523909045002
$ws.Range("O3").Value = 0.3842514532634088
$ws.Range("P3").Value = 0.3842514532634088
$ws.Range("Q3").Value = 1.336566817109333
$ws.Range("R3").Value = 12.029101353984
$ws.Range("S3").Value = 0.3598606730664705
$ws.Range("T3").Value = 0.3598606730664705

# Row 4
$ws.Range("G4").Value = 0.1693853333333334
$ws.Range("H4").Value = 0.508156
$ws.Range("I4").Value = 0.936523909045002
$ws.Range("J4").Value = 0.936523909045002
$ws.Range("M4").Value = 4.974008666666667
$ws.Range("N4").Value = 14.922026
$ws.Range("O4").Value = 0.2422184299659874
$ws.Range("P4").Value = 0.2422184299659874
$ws.Range("Q4").Value = 0.8425241160062225
$ws.Range("R4").Value = 7.582717044056001
$ws.Range("S4").Value = 0.2268433508744896
$ws.Range("T4").Value = 0.2268433508744896

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.1693853333333334
$ws.Range("H5").Value = 0.508156
$ws.Range("I5").Value = 0.936523909045002
$ws.Range("J5").Value = 0.936523909045002
$ws.Range("M5").Value = 2.087648
$ws.Range("N5").Value = 6.262943999999999
$ws.Range("O5").Value = 0.1016618294757629
$ws.Range("P5").Value = 0.1016618294757629
$ws.Range("Q5").Value = 0.3536169523626667
$ws.Range("R5").Value = 3.182552571264
$ws.Range("S5").Value = 0.09520873394130791
$ws.Range("T5").Value = 0.0952087339413079

# Row 6
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.1693853333333334
$ws.Range("H6").Value = 0.508156
$ws.Range("I6").Value = 0.936523909045002
$ws.Range("J6").Value = 0.936523909045002
$ws.Range("M6").Value = 4.848944666666667
$ws.Range("N6").Value = 14.546834
$ws.Range("O6").Value = 0.2361282102347124
$ws.Range("P6").Value = 0.2361282102347124
$ws.Range("Q6").Value = 0.8213401086782224
$ws.Range("R6").Value = 7.392060978104001
$ws.Range("S6").Value = 0.2211397144848129
$ws.Range("T6").Value = 0.2211397144848129

# Row 7
$ws.Range("G7").Value = 0.01148066666666667
$ws.Range("H7").Value = 0.034442
$ws.Range("I7").Value = 0.06347609095499798
$ws.Range("J7").Value = 0.06347609095499798
$ws.Range("M7").Value = 0.7339303333333334
$ws.Range("N7").Value = 2.201791
$ws.Range("O7").Value = 0.03574007706012852
$ws.Range("P7").Value = 0.03574007706012852
$ws.Range("Q7").Value = 0.008426009513555556
$ws.Range("R7").Value = 0.075834085622
$ws.Range("S7").Value = 0.002268640382207355
$ws.Range("T7").Value = 0.002268640382207354

# Row 8
$ws.Range("G8").Value = 0.01148066666666667
$ws.Range("H8").Value = 0.034442
$ws.Range("I8").Value = 0.06347609095499798
$ws.Range("J8").Value = 0.06347609095499798
$ws.Range("O8").Value = 0.3842514532634088
$ws.Range("P8").Value = 0.3842514532634088
$ws.Range("Q8").Value = 0.09059035869866666
$ws.Range("R8").Value = 0.815313228288
$ws.Range("S8").Value = 0.02439078019693829
$ws.Range("T8").Value = 0.02439078019693829

# Row 9
$ws.Range("G9").Value = 0.01148066666666667
$ws.Range("H9").Value = 0.034442
$ws.Range("I9").Value = 0.06347609095499798
$ws.Range("J9").Value = 0.06347609095499798
$ws.Range("M9").Value = 4.974008666666667
$ws.Range("N9").Value = 14.922026
$ws.Range("O9").Value = 0.2422184299659874
$ws.Range("P9").Value = 0.2422184299659874
$ws.Range("Q9").Value = 0.05710493549911112
$ws.Range("R9").Value = 0.513944419492
$ws.Range("S9").Value = 0.01537507909149783
$ws.Range("T9").Value = 0.01537507909149782

# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("G10").Value = 0.01148066666666667
$ws.Range("H10").Value = 0.034442
$ws.Range("I10").Value = 0.06347609095499798
$ws.Range("J10").Value = 0.06347609095499798
$ws.Range("M10").Value = 2.087648
$ws.Range("N10").Value = 6.262943999999999
$ws.Range("O10").Value = 0.1016618294757629
$ws.Range("P10").Value = 0.1016618294757629
$ws.Range("Q10").Value = 0.02396759080533333
$ws.Range("R10").Value = 0.215708317248
$ws.Range("S10").Value = 0.006453095534455023
$ws.Range("T10").Value = 0.006453095534455022

# Row 11
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("G11").Value = 0.01148066666666667
$ws.Range("H11").Value = 0.034442
$ws.Range("I11").Value = 0.06347609095499798
$ws.Range("J11").Value = 0.06347609095499798
$ws.Range("M11").Value = 4.848944666666667
$ws.Range("N11").Value = 14.546834
$ws.Range("O11").Value = 0.2361282102347124
$ws.Range("P11").Value = 0.2361282102347124
$ws.Range("Q11").Value = 0.05566911740311112
$ws.Range("R11").Value = 0.501022056628
$ws.Range("S11").Value = 0.01498849574989949
$ws.Range("T11").Value = 0.01498849574989949
